$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# zh-cn sheet: the "bc2f3d67..." row (row 2) has been handed back and is now
# in sync with en-US. Update its status, target/handback file links and the
# handback datetime.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/496408be7e455685d426d2bb4f31c5f9388c3e51/e2e/bc2f3d67-d464-4b0f-b2a7-526819421d68.md",
    "",
    "",
    "bc2f3d67-d464-4b0f-b2a7-526819421d68.md"
)

$wsZh.Hyperlinks.Add(
    $wsZh.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbd45ae0be2b589de182270cd61a1c8da6640818/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.zh-cn.xlf",
    "",
    "",
    "bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.zh-cn.xlf"
)

$wsZh.Range("H2").Value = "2016-03-12 20:32:08"

# ---------------------------------------------------------------------------
# de-de sheet: same update for the "bc2f3d67..." row.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/496408be7e455685d426d2bb4f31c5f9388c3e51/e2e/bc2f3d67-d464-4b0f-b2a7-526819421d68.md",
    "",
    "",
    "bc2f3d67-d464-4b0f-b2a7-526819421d68.md"
)

$wsDe.Hyperlinks.Add(
    $wsDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/66e474cc1a5ba165e43a96a2a040757a3900033d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.de-de.xlf",
    "",
    "",
    "bc2f3d67-d464-4b0f-b2a7-526819421d68.85ca65e91160a90df44eac3a42a67869e984e0fb.de-de.xlf"
)

$wsDe.Range("H2").Value = "2016-03-12 20:32:14"
